# nmi/comments.xlsx: B1 header now holds the query refresh's date-stamp
# "8/1/2023" (stored as text but shown with a date number format), which
# also renames the linked table's second column to match. Column B is
# narrowed to a fixed custom width (no longer auto "best fit"), and the
# active selection ends up on B1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cell = $ws.Range("B1")
$cell.NumberFormat = "@"
$cell.Value = "8/1/2023"
$cell.NumberFormat = "mm-dd-yy"

$ws.Columns.Item(2).ColumnWidth = 86.3

$ws.Range("B1").Select() | Out-Null
